# XtEHR "Device and DeviceUse" map update.
#
# The EHDSDeviceUse block (rows 15-28 in the old layout) is restructured:
# a set of "header.*" fields is introduced, a few leaf fields are renamed
# ([x] suffixes, "startDate" instead of "implantDate", etc.) and the whole
# block grows from 14 to 24 rows, pushing the trailing MedicalDevice-only
# rows down accordingly. Rows 1-14 (EHDSDevice / MedicalDevice.Product
# section) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(15, 1).Value = 'EHDSDeviceUse.header'
$ws.Cells.Item(15, 2).Value = ''

$ws.Cells.Item(16, 1).Value = 'EHDSDeviceUse.header.subject'
$ws.Cells.Item(16, 2).Value = ''

$ws.Cells.Item(17, 1).Value = 'EHDSDeviceUse.header.identifier'
$ws.Cells.Item(17, 2).Value = ''

$ws.Cells.Item(18, 1).Value = 'EHDSDeviceUse.header.authorship'
$ws.Cells.Item(18, 2).Value = ''

$ws.Cells.Item(19, 1).Value = 'EHDSDeviceUse.header.authorship.author[x]'
$ws.Cells.Item(19, 2).Value = ''

$ws.Cells.Item(20, 1).Value = 'EHDSDeviceUse.header.authorship.datetime'
$ws.Cells.Item(20, 2).Value = ''

$ws.Cells.Item(21, 1).Value = 'EHDSDeviceUse.header.lastUpdate'
$ws.Cells.Item(21, 2).Value = ''

$ws.Cells.Item(22, 1).Value = 'EHDSDeviceUse.header.status'
$ws.Cells.Item(22, 2).Value = ''

$ws.Cells.Item(23, 1).Value = 'EHDSDeviceUse.header.statusReason[x]'
$ws.Cells.Item(23, 2).Value = ''

$ws.Cells.Item(24, 1).Value = 'EHDSDeviceUse.header.language'
$ws.Cells.Item(24, 2).Value = ''

$ws.Cells.Item(25, 1).Value = 'EHDSDeviceUse.header.version'
$ws.Cells.Item(25, 2).Value = ''

$ws.Cells.Item(26, 1).Value = 'EHDSDeviceUse.presentedForm'
$ws.Cells.Item(26, 2).Value = ''

$ws.Cells.Item(27, 1).Value = 'EHDSDeviceUse.status'
$ws.Cells.Item(27, 2).Value = ''

$ws.Cells.Item(28, 1).Value = 'EHDSDeviceUse.startDate'
$ws.Cells.Item(28, 2).Value = 'MedicalDevice.StartDate'

$ws.Cells.Item(29, 1).Value = 'EHDSDeviceUse.endDate'
$ws.Cells.Item(29, 2).Value = 'MedicalDevice.EndDate'

$ws.Cells.Item(30, 1).Value = 'EHDSDeviceUse.device[x]'
$ws.Cells.Item(30, 2).Value = ''

$ws.Cells.Item(31, 1).Value = 'EHDSDeviceUse.bodySite'
$ws.Cells.Item(31, 2).Value = 'MedicalDevice.AnatomicalLocation'

$ws.Cells.Item(32, 1).Value = 'EHDSDeviceUse.note'
$ws.Cells.Item(32, 2).Value = 'MedicalDevice.Comment'

$ws.Cells.Item(33, 1).Value = 'EHDSDeviceUse.recorded'
$ws.Cells.Item(33, 2).Value = ''

# Row 38 ("reason[x]") is written before row 34 ("source[x]") so new
# shared-string entries are registered in the same order as the source
# workbook.
$ws.Cells.Item(38, 1).Value = 'EHDSDeviceUse.reason[x]'
$ws.Cells.Item(38, 2).Value = 'MedicalDevice.Indication::Diagnosis'

$ws.Cells.Item(34, 1).Value = 'EHDSDeviceUse.source[x]'
$ws.Cells.Item(34, 2).Value = ''

$ws.Cells.Item(35, 1).Value = ''
$ws.Cells.Item(35, 2).Value = 'MedicalDevice.ProductDescription'

$ws.Cells.Item(36, 1).Value = ''
$ws.Cells.Item(36, 2).Value = 'MedicalDevice.Location::HealthcareProvider'

$ws.Cells.Item(37, 1).Value = ''
$ws.Cells.Item(37, 2).Value = 'MedicalDevice.HealthProfessional'

# Selection/view moved as the sheet grew (no more frozen/scrolled topLeftCell).
$ws.Range("B33").Select() | Out-Null

# Printer settings were touched in this revision as well.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
